$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.305.40"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.705.93"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "223.96"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "0.5306"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").Value = "0.06578"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").Value = "20.73"
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("D11").Value = "0.07621"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").Value = "4.512"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").Value = "1.708.08"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "1.939.42"
$ws.Range("D15").Value = "0.5768"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").Value = "0.0₅8137"
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("D17").Value = "67.59"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "27.302.86"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "215.39"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("D23").Value = "5.955"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "144.38"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("D26").Value = "1.711"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "0.1200"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("D28").Value = "7.211"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").Value = "16.10"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("D30").Value = "0.05376"
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("D31").Value = "1.287"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "3.466"
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").Value = "1.645"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").Value = "2.863"
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("D36").Value = "2.415"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("D37").Value = "0.9447"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "0.5791"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "0.01624"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "5.761"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8408"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.039.45"
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("D44").Value = "100.91"
$ws.Range("D45").Value = "1.848.34"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("E46").Value = "  +4.93%  "
$ws.Range("D47").Value = "57.64"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "0.4519"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "8.048"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").Value = "0.05228"
